# 3504 PC to be continued
# Applies the changes described by the OOXML diff to the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Environments_OnGoing: add 5 new Variable/Value rows (9-13)
# ---------------------------------------------------------------------
$wsEnv = $wb.Worksheets.Item("Environments_OnGoing")

$wsEnv.Range("A9").Value = "DirectSalesDiscountConsoleByProductBasket"
$wsEnv.Range("B9").Value = "https://proximus--prxitt.lightning.force.com/apex/csdiscounts__DiscountPage?basketId="

$wsEnv.Range("A10").Value = "DirectSalesOrderEnrichmentByProductBasket"
$wsEnv.Range("B10").Value = "https://proximus--prxitt--csoe.visualforce.com/apex/apex/NonCommercialSpecifications?basketId="

$wsEnv.Range("A11").Value = "DirectSalesOrderEnrichmentByOrder"
$wsEnv.Range("B11").Value = "https://proximus--prxitt--csoe.visualforce.com/apex/apex/NonCommercialSpecifications?orderId="

$wsEnv.Range("A12").Value = "PartnersCommunityOpportunityRelated"
$wsEnv.Range("B12").Value = "/s/opportunity/related"

$wsEnv.Range("A13").Value = "PartnersCommunityOpportunitesSufix"
$wsEnv.Range("B13").Value = "/Opportunities"

# New hyperlinks on B9, B11, B10 -- added in this order so the generated
# relationship ids come out as rId5/rId6/rId7 in the same arrangement as
# the target workbook (B9->rId5, B11->rId6, B10->rId7).
$wsEnv.Hyperlinks.Add($wsEnv.Range("B9"), "https://proximus--prxitt.lightning.force.com/apex/csdiscounts__DiscountPage?basketId=")
$wsEnv.Hyperlinks.Add($wsEnv.Range("B11"), "https://proximus--prxitt--csoe.visualforce.com/apex/apex/NonCommercialSpecifications?orderId=")
$wsEnv.Hyperlinks.Add($wsEnv.Range("B10"), "https://proximus--prxitt--csoe.visualforce.com/apex/apex/NonCommercialSpecifications?basketId=")

# Hyperlinks.Add() stamps a freshly-created style onto the cell; restore the
# workbook's existing shared "Hyperlink" cell style so B9:B11 look the same
# as the pre-existing hyperlinked cells (B2, B3, B4, B7) in this sheet.
$wsEnv.Range("B9").Style = "Hyperlink"
$wsEnv.Range("B10").Style = "Hyperlink"
$wsEnv.Range("B11").Style = "Hyperlink"

# Column A got wider to fit the new, longer variable names.
$wsEnv.Columns.Item(1).ColumnWidth = 40.6

# ---------------------------------------------------------------------
# 2) Users_OnGoing: scrub the personal farmer-profile credentials
# ---------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users_OnGoing")

$wsUsers.Range("B12").Value = "soi.testing.crew@gmail.com.farmer"
$wsUsers.Range("B13").Value = "Testing-2020"

$wsUsers.Hyperlinks.Add($wsUsers.Range("B12"), "mailto:soi.testing.crew@gmail.com.farmer")
$wsUsers.Range("B12").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 3) View-state: selections on a few sheets, and D03NonQuotableProducts
#    becomes the active tab (was Users_OnGoing).
# ---------------------------------------------------------------------
$wsEnv.Activate()
$wsEnv.Range("B17").Select()

$wsCompanies = $wb.Worksheets.Item("Companies")
$wsCompanies.Activate()
$wsCompanies.Range("B38").Select()

$wsD03 = $wb.Worksheets.Item("D03NonQuotableProducts")
$wsD03.Activate()
$wsD03.Range("J9").Select()
